# Update cryptocurrency price and volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.145.93"
$ws.Range("E2").Value = "  -1.07%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.677.53"
$ws.Range("E3").Value = "  -0.83%  "

$ws.Range("E4").Value = "  -0.50%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "210.03"
$ws.Range("E5").Value = "  -4.10%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5267"
$ws.Range("E6").Value = "  -4.81%  "

$ws.Range("E7").Value = "  -0.43%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2674"
$ws.Range("E8").Value = "  -1.55%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06279"
$ws.Range("E9").Value = "  -3.12%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.21"
$ws.Range("E10").Value = "  -4.11%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07518"
$ws.Range("E11").Value = "  -1.28%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.688.70"
$ws.Range("E12").Value = "  -0.19%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.489"
$ws.Range("E13").Value = "  -1.62%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5645"
$ws.Range("E14").Value = "  -3.27%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.000008098"
$ws.Range("E15").Value = "  -4.31%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.00"
$ws.Range("E16").Value = "  +1.00%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.187.85"
$ws.Range("E17").Value = "  -1.18%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.006"
$ws.Range("E18").Value = "  -0.42%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.842"
$ws.Range("E19").Value = "  -2.47%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.51"
$ws.Range("E20").Value = "  -4.32%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "188.33"
$ws.Range("E21").Value = "  -1.27%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.187"
$ws.Range("E22").Value = "  -1.16%  "

$ws.Range("E23").Value = "  -0.45%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "148.03"
$ws.Range("E24").Value = "  -1.36%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1256"
$ws.Range("E25").Value = "  -4.20%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.584"
$ws.Range("E26").Value = "  -4.27%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.81"
$ws.Range("E27").Value = "  +0.07%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.06429"
$ws.Range("E28").Value = "  +1.36%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.336"
$ws.Range("E29").Value = "  -5.60%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.281"
$ws.Range("E30").Value = "  -3.69%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.528"
$ws.Range("E31").Value = "  -1.80%  "

$ws.Range("E32").Value = "  -3.23%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.646"
$ws.Range("E33").Value = "  -1.93%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.003"
$ws.Range("E34").Value = "  -4.29%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6081"
$ws.Range("E35").Value = "  -2.74%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.417"
$ws.Range("E36").Value = "  +0.54%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.713"

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.162"
$ws.Range("E38").Value = "  -1.43%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.098.32"
$ws.Range("E39").Value = "  -2.35%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01605"
$ws.Range("E40").Value = "  -2.63%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8662"
$ws.Range("E41").Value = "  -2.27%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.007"
$ws.Range("E42").Value = "  -0.98%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "100.04"
$ws.Range("E43").Value = "  -0.69%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.830.03"
$ws.Range("E44").Value = "  -0.74%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00000000110"
$ws.Range("E45").Value = "  -0.19%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "56.81"
$ws.Range("E46").Value = "  -1.35%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.001"
$ws.Range("E47").Value = "  -0.83%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05273"
$ws.Range("E48").Value = "  -0.20%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.989"
$ws.Range("E49").Value = "  -3.36%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4274"
$ws.Range("E50").Value = "  -0.64%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.941"
$ws.Range("E51").Value = "  -2.45%  "
